$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Flujo de caja")

# Update the monthly remuneration base values that feed the "U14" (total
# monthly personnel cost) calculation further down the cash-flow sheet.
$ws.Range("R14").Value = 120000
$ws.Range("R15").Value = 70000

$ws.Activate()
$ws.Range("R16").Select()

$wb.Save()
